$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header (shared string) renames ---
$ws.Range("C1").Value = "GDP"
$ws.Range("E1").Value = "Budget_Previous_Year"
$ws.Range("F1").Value = "LatinAmerica"
$ws.Range("G1").Value = "Africa"
$ws.Range("H1").Value = "Confessional"
$ws.Range("I1").Value = "Universal"
$ws.Range("AF1").Value = "Donor_Aid_Budget"
$ws.Range("AG1").Value = "Total_Funds"
$ws.Range("AH1").Value = "%_Private_Funds"
$ws.Range("AI1").Value = "%_MAE_Funds"
$ws.Range("AM1").Value = "Delegation"

# --- Column C (GDP) numeric updates, rows 2-108 ---
$colC = @{
    2 = 4403.873446482281
    3 = 2565.769164761619
    4 = 2771.04675450926
    5 = 4043.14643519898
    6 = 2934.187009790061
    7 = 2870.311589353206
    8 = 8947.741473873051
    9 = 5191.140356354663
    10 = 9502.243585046588
    11 = 6128.19547247793
    12 = 4547.50930098406
    13 = 3972.630273980753
    14 = 4729.735976516416
    16 = 4480.799271301591
    17 = 2644.817039097506
    18 = 2839.92516805933
    19 = 4141.976353364595
    20 = 557.6321326261259
    21 = 710.2742021758368
    22 = 473.2998774917226
    23 = 1280.225469721551
    24 = 2898.942214704482
    25 = 665.6274194933962
    26 = 9271.398233246389
    27 = 5555.389721901988
    28 = 6336.709213679884
    29 = 4633.590358399045
    30 = 4355.934938677345
    31 = 1317.890706178356
    32 = 4524.720276132375
    33 = 2635.191563620458
    34 = 2948.84548976845
    35 = 4022.237688257
    36 = 587.1403830380606
    37 = 711.1128122770988
    38 = 1629.435089125503
    39 = 466.0709276378625
    40 = 1263.452411343738
    41 = 2965.153206179127
    42 = 691.8942672110555
    43 = 9477.887185090232
    44 = 5660.517066940175
    45 = 10883.31535948899
    46 = 6711.616186806423
    47 = 4921.848409120176
    48 = 4479.398934239905
    49 = 1385.890384668919
    50 = 4394.543881413723
    51 = 1443.492614888721
    52 = 4142.494431377753
    53 = 2995.45235738661
    54 = 5122.180090208862
    55 = 2999.422762626143
    56 = 2634.85005236495
    57 = 10649.8372874572
    58 = 495.763971160512
    59 = 561.9866500362131
    60 = 701.4459636783288
    61 = 5745.422744292303
    62 = 4699.493713911862
    63 = 1505.810948829135
    64 = 4220.389835426517
    65 = 3087.12349650562
    66 = 5295.682695961288
    67 = 3056.152683606517
    68 = 2632.058233068435
    69 = 10784.63069920566
    70 = 503.3023574516347
    71 = 565.0479699255185
    72 = 720.7128711178943
    73 = 5955.175904294275
    74 = 7449.08671983612
    75 = 5412.131646018807
    76 = 2648.294169302945
    77 = 3137.260298393558
    78 = 555.6050487299622
    79 = 730.3063521039821
    80 = 3125.07948072635
    81 = 515.8271637832048
    82 = 4861.287098802361
    83 = 6301.696269820412
    84 = 1338.716747746975
    85 = 4302.486544228166
    86 = 1579.189101937001
    87 = 7580.275568826287
    88 = 5330.539154475424
    89 = 2703.742092148914
    90 = 3210.869677115934
    91 = 574.4917005765921
    92 = 729.1196658666737
    93 = 3222.05417836739
    94 = 517.8609592583078
    95 = 4944.191641077407
    96 = 6661.86504232374
    97 = 1384.519227335143
    98 = 4308.415595998206
    99 = 1667.171891046301
    100 = 7633.969039669125
    101 = 5176.058803160127
    102 = 3242.636921959078
    103 = 594.9783586927158
    104 = 3212.81539531051
    105 = 5089.61202008711
    106 = 7026.178156858586
    107 = 1431.756130822538
    108 = 4310.567450080138
}
foreach ($row in $colC.Keys) {
    $ws.Cells.Item($row, 3).Value = $colC[$row]
}

# --- Column AL (Visitado flag) updates: 0 -> 1 ---
$alRows = @(9, 25, 27, 42, 44, 60, 61, 72, 73, 79, 83, 92, 96, 106)
foreach ($row in $alRows) {
    $ws.Cells.Item($row, 38).Value = 1
}
